{"js": "// Auto-generated: replace each math-equation cell's text with its new value,\n// matching cells by their current (old) text so we don't depend on row/col\n// layout assumptions. Values are taken from the table in document order and\n// are each unique, so an old->new text map is a safe, unambiguous plan.\nconst replacements = [\n  [\"8+73=\", \"31-27=\"],\n  [\"54+6=\", \"18+60=\"],\n  [\"92-85=\", \"92+6=\"],\n  [\"0+91=\", \"58+21=\"],\n  [\"45+40=\", \"17+33=\"],\n  [\"22+32=\", \"85-33=\"],\n  [\"62-41=\", \"28+11=\"],\n  [\"23+42=\", \"7+24=\"],\n  [\"65+20=\", \"83+13=\"],\n  [\"0+26=\", \"79+2=\"],\n  [\"51-7=\", \"94-67=\"],\n  [\"51+1=\", \"72-11=\"],\n  [\"77-71=\", \"79-60=\"],\n  [\"60-42=\", \"81-10=\"],\n  [\"90-31=\", \"15+61=\"],\n  [\"22-17=\", \"25-14=\"],\n  [\"82-74=\", \"8+79=\"],\n  [\"28+57=\", \"91-84=\"],\n  [\"78+16=\", \"1+92=\"],\n  [\"95-6=\", \"89-3=\"],\n  [\"33+33=\", \"27+63=\"],\n  [\"53+39=\", \"66-29=\"],\n  [\"2+64=\", \"77-20=\"],\n  [\"20-5=\", \"26+37=\"],\n  [\"82-26=\", \"35+6=\"],\n  [\"23+66=\", \"21+14=\"],\n  [\"77-47=\", \"18+57=\"],\n  [\"7+21=\", \"55-3=\"],\n  [\"78-51=\", \"89-37=\"],\n  [\"58-44=\", \"27+69=\"],\n  [\"90-16=\", \"92-64=\"],\n  [\"28+38=\", \"85-61=\"],\n  [\"64-5=\", \"69-19=\"],\n  [\"74-10=\", \"4+78=\"],\n  [\"7-0=\", \"6+69=\"],\n  [\"51+47=\", \"67-65=\"],\n  [\"87-34=\", \"56+14=\"],\n  [\"34+29=\", \"62-31=\"],\n  [\"51-17=\", \"32+36=\"],\n  [\"58-28=\", \"93-57=\"],\n  [\"85-77=\", \"62+19=\"],\n  [\"27+14=\", \"66-40=\"],\n  [\"48-35=\", \"10+26=\"],\n  [\"61-37=\", \"93-50=\"],\n  [\"77-14=\", \"22+72=\"],\n  [\"71-59=\", \"72+18=\"],\n  [\"55-26=\", \"74-3=\"],\n  [\"29+51=\", \"36+46=\"],\n  [\"84-10=\", \"1+16=\"],\n  [\"42-33=\", \"6+56=\"],\n  [\"34+42=\", \"34+56=\"],\n  [\"43-34=\", \"34+10=\"],\n  [\"44+26=\", \"67+0=\"],\n  [\"36+26=\", \"35+16=\"],\n  [\"16+60=\", \"43+17=\"],\n  [\"33-22=\", \"26+14=\"],\n  [\"24-4=\", \"65+25=\"],\n  [\"79-64=\", \"26+36=\"],\n  [\"47-30=\", \"56+40=\"],\n  [\"18+41=\", \"71-32=\"],\n  [\"79-59=\", \"32-7=\"],\n  [\"93-91=\", \"33+16=\"],\n  [\"59-52=\", \"40+4=\"],\n  [\"90-54=\", \"84-44=\"],\n  [\"7+12=\", \"40-32=\"],\n  [\"16+7=\", \"63-26=\"],\n  [\"27+22=\", \"98-74=\"],\n  [\"61-18=\", \"32-26=\"],\n  [\"63-47=\", \"10+6=\"],\n  [\"69-52=\", \"57+40=\"],\n  [\"36+45=\", \"16+81=\"],\n  [\"68+6=\", \"15+33=\"],\n  [\"58-22=\", \"66+26=\"],\n  [\"83+7=\", \"42+22=\"],\n  [\"32+18=\", \"58-31=\"],\n  [\"67-30=\", \"86-78=\"],\n  [\"19+32=\", \"29+35=\"],\n  [\"19+14=\", \"46+51=\"],\n  [\"67-7=\", \"58+7=\"],\n  [\"24+55=\", \"11+87=\"],\n  [\"65+14=\", \"40+23=\"],\n  [\"1+32=\", \"96-12=\"],\n  [\"91-89=\", \"63-55=\"],\n  [\"74-16=\", \"68-48=\"],\n  [\"83-67=\", \"90-63=\"],\n  [\"37-19=\", \"69-53=\"],\n  [\"41-29=\", \"51+29=\"],\n  [\"59+27=\", \"7+89=\"],\n  [\"54-24=\", \"73+9=\"],\n  [\"94-21=\", \"20+56=\"],\n  [\"7+91=\", \"86-14=\"],\n  [\"24+68=\", \"87-49=\"],\n  [\"85-74=\", \"10-7=\"],\n  [\"74+10=\", \"98-49=\"],\n  [\"77-22=\", \"34-29=\"],\n  [\"70-33=\", \"66+18=\"],\n  [\"29+16=\", \"31-2=\"],\n  [\"35+30=\", \"47+33=\"],\n  [\"82+3=\", \"75+6=\"],\n  [\"68-44=\", \"27+9=\"]\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\n// Build quick lookup of old text -> new text.\nconst map = new Map(replacements);\n\nlet applied = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    // eslint-disable-next-line no-await-in-loop\n    await context.sync();\n\n    const current = para.text;\n    if (map.has(current)) {\n      para.insertText(map.get(current), Word.InsertLocation.replace);\n      applied++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Auto-generated: replace each math-equation table cell's text with its new\n# value. We match cells by their current (old) text, rather than assuming a\n# fixed row/column order, since every old value in this table is unique.\n$replacements = @{\n  '8+73=' = '31-27=';\n  '54+6=' = '18+60=';\n  '92-85=' = '92+6=';\n  '0+91=' = '58+21=';\n  '45+40=' = '17+33=';\n  '22+32=' = '85-33=';\n  '62-41=' = '28+11=';\n  '23+42=' = '7+24=';\n  '65+20=' = '83+13=';\n  '0+26=' = '79+2=';\n  '51-7=' = '94-67=';\n  '51+1=' = '72-11=';\n  '77-71=' = '79-60=';\n  '60-42=' = '81-10=';\n  '90-31=' = '15+61=';\n  '22-17=' = '25-14=';\n  '82-74=' = '8+79=';\n  '28+57=' = '91-84=';\n  '78+16=' = '1+92=';\n  '95-6=' = '89-3=';\n  '33+33=' = '27+63=';\n  '53+39=' = '66-29=';\n  '2+64=' = '77-20=';\n  '20-5=' = '26+37=';\n  '82-26=' = '35+6=';\n  '23+66=' = '21+14=';\n  '77-47=' = '18+57=';\n  '7+21=' = '55-3=';\n  '78-51=' = '89-37=';\n  '58-44=' = '27+69=';\n  '90-16=' = '92-64=';\n  '28+38=' = '85-61=';\n  '64-5=' = '69-19=';\n  '74-10=' = '4+78=';\n  '7-0=' = '6+69=';\n  '51+47=' = '67-65=';\n  '87-34=' = '56+14=';\n  '34+29=' = '62-31=';\n  '51-17=' = '32+36=';\n  '58-28=' = '93-57=';\n  '85-77=' = '62+19=';\n  '27+14=' = '66-40=';\n  '48-35=' = '10+26=';\n  '61-37=' = '93-50=';\n  '77-14=' = '22+72=';\n  '71-59=' = '72+18=';\n  '55-26=' = '74-3=';\n  '29+51=' = '36+46=';\n  '84-10=' = '1+16=';\n  '42-33=' = '6+56=';\n  '34+42=' = '34+56=';\n  '43-34=' = '34+10=';\n  '44+26=' = '67+0=';\n  '36+26=' = '35+16=';\n  '16+60=' = '43+17=';\n  '33-22=' = '26+14=';\n  '24-4=' = '65+25=';\n  '79-64=' = '26+36=';\n  '47-30=' = '56+40=';\n  '18+41=' = '71-32=';\n  '79-59=' = '32-7=';\n  '93-91=' = '33+16=';\n  '59-52=' = '40+4=';\n  '90-54=' = '84-44=';\n  '7+12=' = '40-32=';\n  '16+7=' = '63-26=';\n  '27+22=' = '98-74=';\n  '61-18=' = '32-26=';\n  '63-47=' = '10+6=';\n  '69-52=' = '57+40=';\n  '36+45=' = '16+81=';\n  '68+6=' = '15+33=';\n  '58-22=' = '66+26=';\n  '83+7=' = '42+22=';\n  '32+18=' = '58-31=';\n  '67-30=' = '86-78=';\n  '19+32=' = '29+35=';\n  '19+14=' = '46+51=';\n  '67-7=' = '58+7=';\n  '24+55=' = '11+87=';\n  '65+14=' = '40+23=';\n  '1+32=' = '96-12=';\n  '91-89=' = '63-55=';\n  '74-16=' = '68-48=';\n  '83-67=' = '90-63=';\n  '37-19=' = '69-53=';\n  '41-29=' = '51+29=';\n  '59+27=' = '7+89=';\n  '54-24=' = '73+9=';\n  '94-21=' = '20+56=';\n  '7+91=' = '86-14=';\n  '24+68=' = '87-49=';\n  '85-74=' = '10-7=';\n  '74+10=' = '98-49=';\n  '77-22=' = '34-29=';\n  '70-33=' = '66+18=';\n  '29+16=' = '31-2=';\n  '35+30=' = '47+33=';\n  '82+3=' = '75+6=';\n  '68-44=' = '27+9='\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $raw = $cell.Range.Text\n    # Strip the trailing cell-mark / paragraph-mark characters (chr 13, chr 7).\n    $current = $raw.TrimEnd([char]7, [char]13)\n    if ($replacements.ContainsKey($current)) {\n      $cell.Range.Text = $replacements[$current]\n    }\n  }\n}\n"}
